# [Complex] Add SceneEnd & Dodge & Critical
#
# The quest descriptions in column C used literal "\n" sequences to mark
# line breaks. They were reworked to use "@" as the separator instead
# (used by the game's new multi-page dialog / SceneEnd system). Re-setting
# these cell values with the new separator causes Excel to register them
# as new shared strings, which is why the shared string table ends up
# re-ordered (new text appended at the end) after the edit - matching the
# target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value  = "당신은 저희 알선소에서 첫 번째 임무를 맞이하고 계십니다.@그래서 역량을 확인하기 위한 간단한 임무로 구성했어요.@첫 모험에 행운이 가득하기를 빌겠습니다."
$ws.Range("C3").Value  = "최근 늑대들이 농산물을 먹어치우고 있다는 신고가 다발하고 있어요.@저희 알선소에서도 적극 지원하고 있지만 해결은 요원하네요.@부디 손을 거들어주시기를 바라겠습니다."
$ws.Range("C4").Value  = "그거 아세요? 원래 슬라임은 슬라루의 제왕개체였다고 해요.@그래서인지 슬라임의 체액이 보약이라는 소문도 도처에 퍼져있구요.@그딴 거 없는데 말이죠. 괜한 피해자가 나오지 않게 수를 줄여주세요."
$ws.Range("C5").Value  = "아... 결국에는 민간인이 슬라임한테 돌격하는 사건이 발생해버렸어요.@엄청 큰 슬라임을 만나서 '전설의 슬라임이다!' 하고 달려들었다네요.@보고 상으로는 평균 크기 보다 조금 작았습니다."
$ws.Range("C6").Value  = "가도에 엄청 큰 늑대가 발견됬다는 보고가 잇다르고 있어요.@어차피 모험가들이 자주 다니는 길목이라 놔둬도 되기는 하는데…@아무래도 사람을 습격한다고 합니다."
$ws.Range("C7").Value  = "늑대들 중에서는 태양의 사자라고 불리는 개체가 있다고 해요.@사실인지는 모르겠어요. 근데 뭔가 고상해 보이기는 해요. 흔하지만.@걔네들이 최근 근처 상단을 습격했다고 하니까 개체 수 좀 줄여주세요."
$ws.Range("C8").Value  = "늑대들 중에서는 달의 사자라고 불리는 개체가…@이건 이전에도 설명 드렸을테니 각설하고.@이번에는 얘네들이 집단 서리를 했대요. 꽤나 많이 털린 모양인지 농장주 께서 좀 많이 화나셨다고…"
$ws.Range("C9").Value  = "최근 도시 근처 연못에서 크지는 않은데 뭔가 고급진 슬라임이 등장했다고 해요.@마찬가지로 '보약 아니냐?!'라는 소문이 계속해서 퍼진다고 하네요.@슬라임 따위한테 그딴거 없을테니 토벌 해 주세요."
$ws.Range("C10").Value = "최근 늑대무리의 이상현상의 원인이 밝혀졌습니다.@늑대들의 신... 과거의 찬달자... 세계의 구원자…@그 모독적인 껍데기를 처치해 주세요."

# Update the active selection, matching the cursor position recorded in the
# saved workbook.
$ws.Range("I8").Select()
